# Fill in attendance (P/F) for the new day column AM on the "Chamada" sheet,
# and update the active cell selection to AM1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

$ws.Range("AM3").Value = "P"
$ws.Range("AM4").Value = "P"
$ws.Range("AM5").Value = "F"
$ws.Range("AM6").Value = "P"
$ws.Range("AM7").Value = "P"
$ws.Range("AM8").Value = "P"
$ws.Range("AM9").Value = "P"
$ws.Range("AM10").Value = "P"
$ws.Range("AM11").Value = "P"
$ws.Range("AM12").Value = "P"
$ws.Range("AM13").Value = "P"
$ws.Range("AM14").Value = "P"
$ws.Range("AM15").Value = "F"
$ws.Range("AM16").Value = "P"
$ws.Range("AM17").Value = "P"
$ws.Range("AM18").Value = "P"
$ws.Range("AM19").Value = "P"
$ws.Range("AM20").Value = "F"
$ws.Range("AM21").Value = "P"
$ws.Range("AM22").Value = "P"
$ws.Range("AM23").Value = "P"
$ws.Range("AM24").Value = "P"
$ws.Range("AM25").Value = "F"
$ws.Range("AM26").Value = "P"
$ws.Range("AM27").Value = "P"
$ws.Range("AM28").Value = "P"
$ws.Range("AM29").Value = "P"

# Update the active selection to reflect where the author left the cursor.
$ws.Activate()
$ws.Range("AM1").Select()
